# Insert two new Cilantro price-report rows at the top of the weekly
# block (rows 285-286), shifting the existing rows 285:317 down to
# 287:319 (dimension grows from A1:R317 to A1:R319).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("285:286").Insert()

# New row 285 - $/caja 36 atados
$ws.Range("A285").Value = 9
$ws.Range("B285").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C285").Value = "Metropolitana"
$ws.Range("D285").Value = 44449
$ws.Range("E285").Value = 13
$ws.Range("F285").Value = 100112040
$ws.Range("G285").Value = "Cilantro"
$ws.Range("H285").Value = "Sin especificar"
$ws.Range("I285").Value = "Primera"
$ws.Range("J285").Value = 43
$ws.Range("K285").Value = 6000
$ws.Range("L285").Value = 6000
$ws.Range("M285").Value = 6000
$ws.Range("N285").Value = "$/caja 36 atados"
$ws.Range("O285").Value = "Región Metropolitana"
$ws.Range("P285").Value = 167
$ws.Range("Q285").Value = 36
$ws.Range("R285").Value = "Hortaliza"

# New row 286 - $/docena de atados
$ws.Range("A286").Value = 9
$ws.Range("B286").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C286").Value = "Metropolitana"
$ws.Range("D286").Value = 44449
$ws.Range("E286").Value = 13
$ws.Range("F286").Value = 100112040
$ws.Range("G286").Value = "Cilantro"
$ws.Range("H286").Value = "Sin especificar"
$ws.Range("I286").Value = "Primera"
$ws.Range("J286").Value = 124
$ws.Range("K286").Value = 10000
$ws.Range("L286").Value = 12000
$ws.Range("M286").Value = 11000
$ws.Range("N286").Value = "$/docena de atados"
$ws.Range("O286").Value = "Región Metropolitana"
$ws.Range("P286").Value = 3667
$ws.Range("Q286").Value = 3
$ws.Range("R286").Value = "Hortaliza"
